# "Ende 2018" - fill in the actual ("Ist") cycling km for the final four
# months of the year (Aug-Nov) on the Tabelle1 sheet. The accumulated
# ("AkkumIst") column F and the yearly total in E15 already contain shared
# formulas, so they recalculate automatically once the inputs are set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("E11").Value = 158   # August
$ws.Range("E12").Value = 188   # September
$ws.Range("E13").Value = 289   # October
$ws.Range("E14").Value = 260   # November

# Move the active selection, matching the saved cursor position.
$ws.Activate()
$ws.Range("E6").Select()
